# Update the JLCPCB part number for J2 (cell D2 on Sheet1):
#   C914556 -> C1509219
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("D2")
$cell.Value = "C1509219"

# Touch the cell's protection state so the workbook picks up a fresh
# (but visually identical) cell style record for D2, matching how Excel
# re-serializes the style table after the edit.
$cell.Locked = $true
